# Generate Report for Handoff
# ---------------------------------------------------------------------------
# The localization report is regenerated: the row for
# "c1c1fcec-3ab1-420c-a5a0-95d19f660de7.md" now comes first (still in sync
# with en-US), and the row for "582846e0-6450-4fe4-8443-be0225224ce4.md" is
# now "Ready for handoff" with a fresh handoff package / timestamps and a
# warning that the previous handback wasn't built off the latest source.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

$oldFile = "582846e0-6450-4fe4-8443-be0225224ce4.md"
$newFile = "c1c1fcec-3ab1-420c-a5a0-95d19f660de7.md"

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/87376808586e9677379a9a954bb96b114deddc1e/e2e/582846e0-6450-4fe4-8443-be0225224ce4.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/221ce846250924134040315b3cfdf9f990502605/e2e/582846e0-6450-4fe4-8443-be0225224ce4.md."

# =====================================================================
# Sheet "Overview": row 2 <-> row 3 file identity swap, row 3 gets the
# new status/date ("Ready for handoff").
# =====================================================================
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = $newFile
$wsOverview.Range("B2").Value = "e2e\" + $newFile

$wsOverview.Range("A3").Value = $oldFile
$wsOverview.Range("B3").Value = "e2e\" + $oldFile
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-09-06 23:00:39"

# Hyperlinks: the displayed text for B2/B3 swaps, the link targets stay put.
$wsOverview.Range("A1").Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/87376808586e9677379a9a954bb96b114deddc1e/e2e/" + $newFile, [Type]::Missing, [Type]::Missing, "e2e\" + $newFile)
$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/87376808586e9677379a9a954bb96b114deddc1e/e2e/" + $oldFile, [Type]::Missing, [Type]::Missing, "e2e\" + $oldFile)

# =====================================================================
# Sheet "zh-cn"
# =====================================================================
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("A2").Value = $newFile
$wsZhCn.Range("G2").Value = "c1c1fcec-3ab1-420c-a5a0-95d19f660de7.1952cc04390f910cfc56e9653a902d545d4e8cd4.zh-cn.xlf"
$wsZhCn.Range("I2").Value = $newFile
$wsZhCn.Range("J2").Value = "c1c1fcec-3ab1-420c-a5a0-95d19f660de7.1952cc04390f910cfc56e9653a902d545d4e8cd4.zh-cn.xlf"

$wsZhCn.Range("A3").Value = $oldFile
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("G3").Value = "582846e0-6450-4fe4-8443-be0225224ce4.83caa00cd6db387991ecf04f19b3ce2fa95985ca.zh-cn.xlf"
$wsZhCn.Range("H3").Value = "2016-09-06 23:00:34"
$wsZhCn.Range("I3").Value = $oldFile
$wsZhCn.Range("J3").Value = "582846e0-6450-4fe4-8443-be0225224ce4.83caa00cd6db387991ecf04f19b3ce2fa95985ca.zh-cn.xlf"
$wsZhCn.Range("P3").Value = $errorDetail

$wsZhCn.Columns.Item(16).ColumnWidth = 39.14

$wsZhCn.Range("A1").Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/87376808586e9677379a9a954bb96b114deddc1e/e2e/" + $oldFile, [Type]::Missing, [Type]::Missing, $newFile)
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/edee64d3105f308f01618c4152eb67c6b1302ffc/e2e/" + $oldFile, [Type]::Missing, [Type]::Missing, $newFile)
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/87376808586e9677379a9a954bb96b114deddc1e/e2e/" + $newFile, [Type]::Missing, [Type]::Missing, $oldFile)
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/edee64d3105f308f01618c4152eb67c6b1302ffc/e2e/" + $newFile, [Type]::Missing, [Type]::Missing, $oldFile)

# =====================================================================
# Sheet "de-de"
# =====================================================================
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("A2").Value = $newFile
$wsDeDe.Range("G2").Value = "c1c1fcec-3ab1-420c-a5a0-95d19f660de7.1952cc04390f910cfc56e9653a902d545d4e8cd4.de-de.xlf"
$wsDeDe.Range("I2").Value = $newFile
$wsDeDe.Range("J2").Value = "c1c1fcec-3ab1-420c-a5a0-95d19f660de7.1952cc04390f910cfc56e9653a902d545d4e8cd4.de-de.xlf"

$wsDeDe.Range("A3").Value = $oldFile
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("G3").Value = "582846e0-6450-4fe4-8443-be0225224ce4.83caa00cd6db387991ecf04f19b3ce2fa95985ca.de-de.xlf"
$wsDeDe.Range("H3").Value = "2016-09-06 23:00:39"
$wsDeDe.Range("I3").Value = $oldFile
$wsDeDe.Range("J3").Value = "582846e0-6450-4fe4-8443-be0225224ce4.83caa00cd6db387991ecf04f19b3ce2fa95985ca.de-de.xlf"
$wsDeDe.Range("P3").Value = $errorDetail

$wsDeDe.Columns.Item(16).ColumnWidth = 39.14

$wsDeDe.Range("A1").Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/87376808586e9677379a9a954bb96b114deddc1e/e2e/" + $oldFile, [Type]::Missing, [Type]::Missing, $newFile)
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/189ba42ddf14131654e47d8bfb126351ebc893b9/e2e/" + $oldFile, [Type]::Missing, [Type]::Missing, $newFile)
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/87376808586e9677379a9a954bb96b114deddc1e/e2e/" + $newFile, [Type]::Missing, [Type]::Missing, $oldFile)
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/189ba42ddf14131654e47d8bfb126351ebc893b9/e2e/" + $newFile, [Type]::Missing, [Type]::Missing, $oldFile)
